# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# --- Cell edits on Tabelle1 row 24: "x" (text) -> 1 (number) -------------
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("K24").Value = 1

# --- View state: scroll the frozen pane and move the active selection ----
$win = $excel.ActiveWindow
[void]$ws.Range("K23").Select()
$win.ScrollRow = 15
$win.ScrollColumn = 5

# --- Resize the workbook window (bookViews/workbookView) -----------------
$win.Width = 13410
$win.Height = 7515
